$d = $word.ActiveDocument

# Phase 1: replace each original value with a unique placeholder token
# to avoid chained-replacement collisions when a new value equals
# another entry old value (e.g. 15+6=21 -> 38+35=73 -> 13+7=20).
$d.Content.Find.Execute("82+8=90", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH000@@", 2) | Out-Null
$d.Content.Find.Execute("16-12=4", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH001@@", 2) | Out-Null
$d.Content.Find.Execute("36+41=77", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH002@@", 2) | Out-Null
$d.Content.Find.Execute("98-52=46", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH003@@", 2) | Out-Null
$d.Content.Find.Execute("60-36=24", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH004@@", 2) | Out-Null
$d.Content.Find.Execute("61+10=71", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH005@@", 2) | Out-Null
$d.Content.Find.Execute("44+53=97", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH006@@", 2) | Out-Null
$d.Content.Find.Execute("85-8=77", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH007@@", 2) | Out-Null
$d.Content.Find.Execute("28-21=7", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH008@@", 2) | Out-Null
$d.Content.Find.Execute("8-1=7", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH009@@", 2) | Out-Null
$d.Content.Find.Execute("66-33=33", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH010@@", 2) | Out-Null
$d.Content.Find.Execute("10-4=6", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH011@@", 2) | Out-Null
$d.Content.Find.Execute("30+51=81", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH012@@", 2) | Out-Null
$d.Content.Find.Execute("74-47=27", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH013@@", 2) | Out-Null
$d.Content.Find.Execute("24+17=41", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH014@@", 2) | Out-Null
$d.Content.Find.Execute("63-31=32", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH015@@", 2) | Out-Null
$d.Content.Find.Execute("65-28=37", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH016@@", 2) | Out-Null
$d.Content.Find.Execute("89-27=62", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH017@@", 2) | Out-Null
$d.Content.Find.Execute("87-22=65", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH018@@", 2) | Out-Null
$d.Content.Find.Execute("44+26=70", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH019@@", 2) | Out-Null
$d.Content.Find.Execute("67-17=50", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH020@@", 2) | Out-Null
$d.Content.Find.Execute("81-50=31", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH021@@", 2) | Out-Null
$d.Content.Find.Execute("67-28=39", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH022@@", 2) | Out-Null
$d.Content.Find.Execute("65-24=41", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH023@@", 2) | Out-Null
$d.Content.Find.Execute("11+29=40", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH024@@", 2) | Out-Null
$d.Content.Find.Execute("72-72=0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH025@@", 2) | Out-Null
$d.Content.Find.Execute("88-86=2", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH026@@", 2) | Out-Null
$d.Content.Find.Execute("5+7=12", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH027@@", 2) | Out-Null
$d.Content.Find.Execute("88-28=60", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH028@@", 2) | Out-Null
$d.Content.Find.Execute("2+97=99", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH029@@", 2) | Out-Null
$d.Content.Find.Execute("62+6=68", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH030@@", 2) | Out-Null
$d.Content.Find.Execute("25+54=79", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH031@@", 2) | Out-Null
$d.Content.Find.Execute("40-40=0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH032@@", 2) | Out-Null
$d.Content.Find.Execute("8+9=17", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH033@@", 2) | Out-Null
$d.Content.Find.Execute("32+61=93", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH034@@", 2) | Out-Null
$d.Content.Find.Execute("66-18=48", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH035@@", 2) | Out-Null
$d.Content.Find.Execute("69+3=72", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH036@@", 2) | Out-Null
$d.Content.Find.Execute("47-0=47", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH037@@", 2) | Out-Null
$d.Content.Find.Execute("74-7=67", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH038@@", 2) | Out-Null
$d.Content.Find.Execute("25+34=59", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH039@@", 2) | Out-Null
$d.Content.Find.Execute("46-15=31", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH040@@", 2) | Out-Null
$d.Content.Find.Execute("93-83=10", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH041@@", 2) | Out-Null
$d.Content.Find.Execute("54-31=23", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH042@@", 2) | Out-Null
$d.Content.Find.Execute("79-58=21", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH043@@", 2) | Out-Null
$d.Content.Find.Execute("0+17=17", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH044@@", 2) | Out-Null
$d.Content.Find.Execute("82-2=80", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH045@@", 2) | Out-Null
$d.Content.Find.Execute("72-69=3", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH046@@", 2) | Out-Null
$d.Content.Find.Execute("34+55=89", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH047@@", 2) | Out-Null
$d.Content.Find.Execute("36-22=14", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH048@@", 2) | Out-Null
$d.Content.Find.Execute("22-21=1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH049@@", 2) | Out-Null
$d.Content.Find.Execute("7+4=11", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH050@@", 2) | Out-Null
$d.Content.Find.Execute("1+39=40", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH051@@", 2) | Out-Null
$d.Content.Find.Execute("3+12=15", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH052@@", 2) | Out-Null
$d.Content.Find.Execute("88+3=91", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH053@@", 2) | Out-Null
$d.Content.Find.Execute("15+6=21", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH054@@", 2) | Out-Null
$d.Content.Find.Execute("9+63=72", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH055@@", 2) | Out-Null
$d.Content.Find.Execute("28+26=54", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH056@@", 2) | Out-Null
$d.Content.Find.Execute("68-42=26", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH057@@", 2) | Out-Null
$d.Content.Find.Execute("14+61=75", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH058@@", 2) | Out-Null
$d.Content.Find.Execute("21+65=86", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH059@@", 2) | Out-Null
$d.Content.Find.Execute("95-11=84", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH060@@", 2) | Out-Null
$d.Content.Find.Execute("24+13=37", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH061@@", 2) | Out-Null
$d.Content.Find.Execute("77-28=49", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH062@@", 2) | Out-Null
$d.Content.Find.Execute("92-68=24", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH063@@", 2) | Out-Null
$d.Content.Find.Execute("22+11=33", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH064@@", 2) | Out-Null
$d.Content.Find.Execute("36-13=23", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH065@@", 2) | Out-Null
$d.Content.Find.Execute("41+35=76", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH066@@", 2) | Out-Null
$d.Content.Find.Execute("51-30=21", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH067@@", 2) | Out-Null
$d.Content.Find.Execute("37+24=61", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH068@@", 2) | Out-Null
$d.Content.Find.Execute("57+9=66", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH069@@", 2) | Out-Null
$d.Content.Find.Execute("95-64=31", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH070@@", 2) | Out-Null
$d.Content.Find.Execute("81-35=46", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH071@@", 2) | Out-Null
$d.Content.Find.Execute("16+46=62", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH072@@", 2) | Out-Null
$d.Content.Find.Execute("20+5=25", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH073@@", 2) | Out-Null
$d.Content.Find.Execute("57+7=64", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH074@@", 2) | Out-Null
$d.Content.Find.Execute("53-46=7", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH075@@", 2) | Out-Null
$d.Content.Find.Execute("44-27=17", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH076@@", 2) | Out-Null
$d.Content.Find.Execute("38+22=60", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH077@@", 2) | Out-Null
$d.Content.Find.Execute("11-10=1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH078@@", 2) | Out-Null
$d.Content.Find.Execute("3+67=70", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH079@@", 2) | Out-Null
$d.Content.Find.Execute("59-45=14", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH080@@", 2) | Out-Null
$d.Content.Find.Execute("32+18=50", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH081@@", 2) | Out-Null
$d.Content.Find.Execute("13+3=16", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH082@@", 2) | Out-Null
$d.Content.Find.Execute("3+3=6", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH083@@", 2) | Out-Null
$d.Content.Find.Execute("96+2=98", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH084@@", 2) | Out-Null
$d.Content.Find.Execute("1+41=42", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH085@@", 2) | Out-Null
$d.Content.Find.Execute("92-66=26", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH086@@", 2) | Out-Null
$d.Content.Find.Execute("73-12=61", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH087@@", 2) | Out-Null
$d.Content.Find.Execute("68-28=40", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH088@@", 2) | Out-Null
$d.Content.Find.Execute("16+4=20", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH089@@", 2) | Out-Null
$d.Content.Find.Execute("89-32=57", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH090@@", 2) | Out-Null
$d.Content.Find.Execute("76-73=3", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH091@@", 2) | Out-Null
$d.Content.Find.Execute("98-74=24", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH092@@", 2) | Out-Null
$d.Content.Find.Execute("49+11=60", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH093@@", 2) | Out-Null
$d.Content.Find.Execute("75-36=39", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH094@@", 2) | Out-Null
$d.Content.Find.Execute("76+18=94", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH095@@", 2) | Out-Null
$d.Content.Find.Execute("4+86=90", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH096@@", 2) | Out-Null
$d.Content.Find.Execute("52-13=39", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH097@@", 2) | Out-Null
$d.Content.Find.Execute("1+30=31", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH098@@", 2) | Out-Null
$d.Content.Find.Execute("38+35=73", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH099@@", 2) | Out-Null

# Phase 2: replace each placeholder token with the final new value
$d.Content.Find.Execute("@@PH000@@", $true, $false, $false, $false, $false, $true, 1, $false, "98-50=48", 2) | Out-Null
$d.Content.Find.Execute("@@PH001@@", $true, $false, $false, $false, $false, $true, 1, $false, "58+24=82", 2) | Out-Null
$d.Content.Find.Execute("@@PH002@@", $true, $false, $false, $false, $false, $true, 1, $false, "7+30=37", 2) | Out-Null
$d.Content.Find.Execute("@@PH003@@", $true, $false, $false, $false, $false, $true, 1, $false, "18+29=47", 2) | Out-Null
$d.Content.Find.Execute("@@PH004@@", $true, $false, $false, $false, $false, $true, 1, $false, "67-24=43", 2) | Out-Null
$d.Content.Find.Execute("@@PH005@@", $true, $false, $false, $false, $false, $true, 1, $false, "42+8=50", 2) | Out-Null
$d.Content.Find.Execute("@@PH006@@", $true, $false, $false, $false, $false, $true, 1, $false, "23+34=57", 2) | Out-Null
$d.Content.Find.Execute("@@PH007@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+41=58", 2) | Out-Null
$d.Content.Find.Execute("@@PH008@@", $true, $false, $false, $false, $false, $true, 1, $false, "92-9=83", 2) | Out-Null
$d.Content.Find.Execute("@@PH009@@", $true, $false, $false, $false, $false, $true, 1, $false, "30+63=93", 2) | Out-Null
$d.Content.Find.Execute("@@PH010@@", $true, $false, $false, $false, $false, $true, 1, $false, "91-64=27", 2) | Out-Null
$d.Content.Find.Execute("@@PH011@@", $true, $false, $false, $false, $false, $true, 1, $false, "0+97=97", 2) | Out-Null
$d.Content.Find.Execute("@@PH012@@", $true, $false, $false, $false, $false, $true, 1, $false, "25-6=19", 2) | Out-Null
$d.Content.Find.Execute("@@PH013@@", $true, $false, $false, $false, $false, $true, 1, $false, "77-52=25", 2) | Out-Null
$d.Content.Find.Execute("@@PH014@@", $true, $false, $false, $false, $false, $true, 1, $false, "53-26=27", 2) | Out-Null
$d.Content.Find.Execute("@@PH015@@", $true, $false, $false, $false, $false, $true, 1, $false, "4+77=81", 2) | Out-Null
$d.Content.Find.Execute("@@PH016@@", $true, $false, $false, $false, $false, $true, 1, $false, "91-79=12", 2) | Out-Null
$d.Content.Find.Execute("@@PH017@@", $true, $false, $false, $false, $false, $true, 1, $false, "11+23=34", 2) | Out-Null
$d.Content.Find.Execute("@@PH018@@", $true, $false, $false, $false, $false, $true, 1, $false, "53-40=13", 2) | Out-Null
$d.Content.Find.Execute("@@PH019@@", $true, $false, $false, $false, $false, $true, 1, $false, "44+39=83", 2) | Out-Null
$d.Content.Find.Execute("@@PH020@@", $true, $false, $false, $false, $false, $true, 1, $false, "19-14=5", 2) | Out-Null
$d.Content.Find.Execute("@@PH021@@", $true, $false, $false, $false, $false, $true, 1, $false, "31+24=55", 2) | Out-Null
$d.Content.Find.Execute("@@PH022@@", $true, $false, $false, $false, $false, $true, 1, $false, "66-20=46", 2) | Out-Null
$d.Content.Find.Execute("@@PH023@@", $true, $false, $false, $false, $false, $true, 1, $false, "45+6=51", 2) | Out-Null
$d.Content.Find.Execute("@@PH024@@", $true, $false, $false, $false, $false, $true, 1, $false, "66+32=98", 2) | Out-Null
$d.Content.Find.Execute("@@PH025@@", $true, $false, $false, $false, $false, $true, 1, $false, "9+81=90", 2) | Out-Null
$d.Content.Find.Execute("@@PH026@@", $true, $false, $false, $false, $false, $true, 1, $false, "4+48=52", 2) | Out-Null
$d.Content.Find.Execute("@@PH027@@", $true, $false, $false, $false, $false, $true, 1, $false, "13+60=73", 2) | Out-Null
$d.Content.Find.Execute("@@PH028@@", $true, $false, $false, $false, $false, $true, 1, $false, "28+47=75", 2) | Out-Null
$d.Content.Find.Execute("@@PH029@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+8=25", 2) | Out-Null
$d.Content.Find.Execute("@@PH030@@", $true, $false, $false, $false, $false, $true, 1, $false, "14-11=3", 2) | Out-Null
$d.Content.Find.Execute("@@PH031@@", $true, $false, $false, $false, $false, $true, 1, $false, "51-6=45", 2) | Out-Null
$d.Content.Find.Execute("@@PH032@@", $true, $false, $false, $false, $false, $true, 1, $false, "11+2=13", 2) | Out-Null
$d.Content.Find.Execute("@@PH033@@", $true, $false, $false, $false, $false, $true, 1, $false, "8+89=97", 2) | Out-Null
$d.Content.Find.Execute("@@PH034@@", $true, $false, $false, $false, $false, $true, 1, $false, "14+56=70", 2) | Out-Null
$d.Content.Find.Execute("@@PH035@@", $true, $false, $false, $false, $false, $true, 1, $false, "91+3=94", 2) | Out-Null
$d.Content.Find.Execute("@@PH036@@", $true, $false, $false, $false, $false, $true, 1, $false, "58-18=40", 2) | Out-Null
$d.Content.Find.Execute("@@PH037@@", $true, $false, $false, $false, $false, $true, 1, $false, "66-36=30", 2) | Out-Null
$d.Content.Find.Execute("@@PH038@@", $true, $false, $false, $false, $false, $true, 1, $false, "8+87=95", 2) | Out-Null
$d.Content.Find.Execute("@@PH039@@", $true, $false, $false, $false, $false, $true, 1, $false, "88-45=43", 2) | Out-Null
$d.Content.Find.Execute("@@PH040@@", $true, $false, $false, $false, $false, $true, 1, $false, "86-60=26", 2) | Out-Null
$d.Content.Find.Execute("@@PH041@@", $true, $false, $false, $false, $false, $true, 1, $false, "20+50=70", 2) | Out-Null
$d.Content.Find.Execute("@@PH042@@", $true, $false, $false, $false, $false, $true, 1, $false, "32+6=38", 2) | Out-Null
$d.Content.Find.Execute("@@PH043@@", $true, $false, $false, $false, $false, $true, 1, $false, "73+11=84", 2) | Out-Null
$d.Content.Find.Execute("@@PH044@@", $true, $false, $false, $false, $false, $true, 1, $false, "92-22=70", 2) | Out-Null
$d.Content.Find.Execute("@@PH045@@", $true, $false, $false, $false, $false, $true, 1, $false, "59-57=2", 2) | Out-Null
$d.Content.Find.Execute("@@PH046@@", $true, $false, $false, $false, $false, $true, 1, $false, "53-42=11", 2) | Out-Null
$d.Content.Find.Execute("@@PH047@@", $true, $false, $false, $false, $false, $true, 1, $false, "16+15=31", 2) | Out-Null
$d.Content.Find.Execute("@@PH048@@", $true, $false, $false, $false, $false, $true, 1, $false, "2+69=71", 2) | Out-Null
$d.Content.Find.Execute("@@PH049@@", $true, $false, $false, $false, $false, $true, 1, $false, "18+42=60", 2) | Out-Null
$d.Content.Find.Execute("@@PH050@@", $true, $false, $false, $false, $false, $true, 1, $false, "88-30=58", 2) | Out-Null
$d.Content.Find.Execute("@@PH051@@", $true, $false, $false, $false, $false, $true, 1, $false, "50-13=37", 2) | Out-Null
$d.Content.Find.Execute("@@PH052@@", $true, $false, $false, $false, $false, $true, 1, $false, "85-83=2", 2) | Out-Null
$d.Content.Find.Execute("@@PH053@@", $true, $false, $false, $false, $false, $true, 1, $false, "36+11=47", 2) | Out-Null
$d.Content.Find.Execute("@@PH054@@", $true, $false, $false, $false, $false, $true, 1, $false, "38+35=73", 2) | Out-Null
$d.Content.Find.Execute("@@PH055@@", $true, $false, $false, $false, $false, $true, 1, $false, "53+18=71", 2) | Out-Null
$d.Content.Find.Execute("@@PH056@@", $true, $false, $false, $false, $false, $true, 1, $false, "23+8=31", 2) | Out-Null
$d.Content.Find.Execute("@@PH057@@", $true, $false, $false, $false, $false, $true, 1, $false, "61+38=99", 2) | Out-Null
$d.Content.Find.Execute("@@PH058@@", $true, $false, $false, $false, $false, $true, 1, $false, "47+5=52", 2) | Out-Null
$d.Content.Find.Execute("@@PH059@@", $true, $false, $false, $false, $false, $true, 1, $false, "30-26=4", 2) | Out-Null
$d.Content.Find.Execute("@@PH060@@", $true, $false, $false, $false, $false, $true, 1, $false, "46-37=9", 2) | Out-Null
$d.Content.Find.Execute("@@PH061@@", $true, $false, $false, $false, $false, $true, 1, $false, "87-35=52", 2) | Out-Null
$d.Content.Find.Execute("@@PH062@@", $true, $false, $false, $false, $false, $true, 1, $false, "92-77=15", 2) | Out-Null
$d.Content.Find.Execute("@@PH063@@", $true, $false, $false, $false, $false, $true, 1, $false, "29+64=93", 2) | Out-Null
$d.Content.Find.Execute("@@PH064@@", $true, $false, $false, $false, $false, $true, 1, $false, "17-11=6", 2) | Out-Null
$d.Content.Find.Execute("@@PH065@@", $true, $false, $false, $false, $false, $true, 1, $false, "10+21=31", 2) | Out-Null
$d.Content.Find.Execute("@@PH066@@", $true, $false, $false, $false, $false, $true, 1, $false, "15+80=95", 2) | Out-Null
$d.Content.Find.Execute("@@PH067@@", $true, $false, $false, $false, $false, $true, 1, $false, "76+0=76", 2) | Out-Null
$d.Content.Find.Execute("@@PH068@@", $true, $false, $false, $false, $false, $true, 1, $false, "54+10=64", 2) | Out-Null
$d.Content.Find.Execute("@@PH069@@", $true, $false, $false, $false, $false, $true, 1, $false, "33+13=46", 2) | Out-Null
$d.Content.Find.Execute("@@PH070@@", $true, $false, $false, $false, $false, $true, 1, $false, "86-55=31", 2) | Out-Null
$d.Content.Find.Execute("@@PH071@@", $true, $false, $false, $false, $false, $true, 1, $false, "12+86=98", 2) | Out-Null
$d.Content.Find.Execute("@@PH072@@", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=0", 2) | Out-Null
$d.Content.Find.Execute("@@PH073@@", $true, $false, $false, $false, $false, $true, 1, $false, "39+50=89", 2) | Out-Null
$d.Content.Find.Execute("@@PH074@@", $true, $false, $false, $false, $false, $true, 1, $false, "4-4=0", 2) | Out-Null
$d.Content.Find.Execute("@@PH075@@", $true, $false, $false, $false, $false, $true, 1, $false, "54-33=21", 2) | Out-Null
$d.Content.Find.Execute("@@PH076@@", $true, $false, $false, $false, $false, $true, 1, $false, "89-74=15", 2) | Out-Null
$d.Content.Find.Execute("@@PH077@@", $true, $false, $false, $false, $false, $true, 1, $false, "98-98=0", 2) | Out-Null
$d.Content.Find.Execute("@@PH078@@", $true, $false, $false, $false, $false, $true, 1, $false, "71+13=84", 2) | Out-Null
$d.Content.Find.Execute("@@PH079@@", $true, $false, $false, $false, $false, $true, 1, $false, "77+7=84", 2) | Out-Null
$d.Content.Find.Execute("@@PH080@@", $true, $false, $false, $false, $false, $true, 1, $false, "25-17=8", 2) | Out-Null
$d.Content.Find.Execute("@@PH081@@", $true, $false, $false, $false, $false, $true, 1, $false, "71-21=50", 2) | Out-Null
$d.Content.Find.Execute("@@PH082@@", $true, $false, $false, $false, $false, $true, 1, $false, "25+3=28", 2) | Out-Null
$d.Content.Find.Execute("@@PH083@@", $true, $false, $false, $false, $false, $true, 1, $false, "30-13=17", 2) | Out-Null
$d.Content.Find.Execute("@@PH084@@", $true, $false, $false, $false, $false, $true, 1, $false, "2+82=84", 2) | Out-Null
$d.Content.Find.Execute("@@PH085@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+60=77", 2) | Out-Null
$d.Content.Find.Execute("@@PH086@@", $true, $false, $false, $false, $false, $true, 1, $false, "69-49=20", 2) | Out-Null
$d.Content.Find.Execute("@@PH087@@", $true, $false, $false, $false, $false, $true, 1, $false, "41-28=13", 2) | Out-Null
$d.Content.Find.Execute("@@PH088@@", $true, $false, $false, $false, $false, $true, 1, $false, "29-10=19", 2) | Out-Null
$d.Content.Find.Execute("@@PH089@@", $true, $false, $false, $false, $false, $true, 1, $false, "75-16=59", 2) | Out-Null
$d.Content.Find.Execute("@@PH090@@", $true, $false, $false, $false, $false, $true, 1, $false, "97-0=97", 2) | Out-Null
$d.Content.Find.Execute("@@PH091@@", $true, $false, $false, $false, $false, $true, 1, $false, "35+55=90", 2) | Out-Null
$d.Content.Find.Execute("@@PH092@@", $true, $false, $false, $false, $false, $true, 1, $false, "66+18=84", 2) | Out-Null
$d.Content.Find.Execute("@@PH093@@", $true, $false, $false, $false, $false, $true, 1, $false, "86-24=62", 2) | Out-Null
$d.Content.Find.Execute("@@PH094@@", $true, $false, $false, $false, $false, $true, 1, $false, "30+52=82", 2) | Out-Null
$d.Content.Find.Execute("@@PH095@@", $true, $false, $false, $false, $false, $true, 1, $false, "86-17=69", 2) | Out-Null
$d.Content.Find.Execute("@@PH096@@", $true, $false, $false, $false, $false, $true, 1, $false, "36-6=30", 2) | Out-Null
$d.Content.Find.Execute("@@PH097@@", $true, $false, $false, $false, $false, $true, 1, $false, "93-89=4", 2) | Out-Null
$d.Content.Find.Execute("@@PH098@@", $true, $false, $false, $false, $false, $true, 1, $false, "42+26=68", 2) | Out-Null
$d.Content.Find.Execute("@@PH099@@", $true, $false, $false, $false, $false, $true, 1, $false, "13+7=20", 2) | Out-Null

Write-Output "Replacements complete"
